# Update the "median_training" sheet with refreshed median/count values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = 58
$ws.Range("B13").Value = 66.83333333333333
$ws.Range("B15").Value = 66.83333333333333
$ws.Range("B16").Value = 16
$ws.Range("B20").Value = 11
$ws.Range("B28").Value = 302
$ws.Range("B29").Value = 329.5
$ws.Range("B30").Value = 551
$ws.Range("B31").Value = 329.5
$ws.Range("B32").Value = 338
$ws.Range("B33").Value = 363.5
$ws.Range("B34").Value = 595
$ws.Range("B35").Value = 363.5
$ws.Range("B36").Value = 79.48999999999999
$ws.Range("B37").Value = 77.545
$ws.Range("B38").Value = 88.43000000000001
$ws.Range("B39").Value = 77.545
